# Update the FHIR StructureDefinition export workbook:
#  - bump the IG version & publish date
#  - replace the Publisher/Contact metadata rows with Publisher/Jurisdiction info
#  - remove the stray duplicated "Contact" row
#  - refresh the Short/Definition text for the root Extension element

$wb = $excel.ActiveWorkbook

# ---- "Metadata" sheet (first tab) ----
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refresh the publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> Alvearie Team
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely,
# shifting all subsequent rows up by one.
$meta.Rows.Item(11).Delete()

# ---- "Elements" sheet (second tab) ----
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element - update its Short (K) and Definition (L) text.
$elements.Range("K2").Value = "Status Code"
$elements.Range("L2").Value = "Customer-specific patient status codes"
